$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column X ("Last Cr") and Column Y ("Last Date") need updated values so the
# last-recorded chromium result/date reflects an actual detected result rather
# than the non-detect default. Values are forced as text to match the original
# inline-string cell typing used throughout this table.
$ws.Range("X2").NumberFormat = "@"
$ws.Range("X2").Value = "0.801"
$ws.Range("Y2").NumberFormat = "@"
$ws.Range("Y2").Value = "2013-07-22"

$ws.Range("X8").NumberFormat = "@"
$ws.Range("X8").Value = "2.11"
$ws.Range("Y8").NumberFormat = "@"
$ws.Range("Y8").Value = "2015-10-26"

$ws.Range("X9").NumberFormat = "@"
$ws.Range("X9").Value = "3.2"
$ws.Range("Y9").NumberFormat = "@"
$ws.Range("Y9").Value = "2010-10-18"

$ws.Range("X10").NumberFormat = "@"
$ws.Range("X10").Value = "2.23"
$ws.Range("Y10").NumberFormat = "@"
$ws.Range("Y10").Value = "2011-10-26"

$ws.Range("X11").NumberFormat = "@"
$ws.Range("X11").Value = "0.92"
$ws.Range("Y11").NumberFormat = "@"
$ws.Range("Y11").Value = "2018-10-22"

$ws.Range("X13").NumberFormat = "@"
$ws.Range("X13").Value = "0.161"
$ws.Range("Y13").NumberFormat = "@"
$ws.Range("Y13").Value = "2013-08-15"

$ws.Range("X14").NumberFormat = "@"
$ws.Range("X14").Value = "2.34"
$ws.Range("Y14").NumberFormat = "@"
$ws.Range("Y14").Value = "2016-08-24"

$ws.Range("X15").NumberFormat = "@"
$ws.Range("X15").Value = "3.31"
$ws.Range("Y15").NumberFormat = "@"
$ws.Range("Y15").Value = "2018-09-06"

$ws.Range("X16").NumberFormat = "@"
$ws.Range("X16").Value = "3.02"
$ws.Range("Y16").NumberFormat = "@"
$ws.Range("Y16").Value = "2016-08-31"

$ws.Range("X17").NumberFormat = "@"
$ws.Range("X17").Value = "2.09"
$ws.Range("Y17").NumberFormat = "@"
$ws.Range("Y17").Value = "2013-06-03"

$ws.Range("X18").NumberFormat = "@"
$ws.Range("X18").Value = "2.85"
$ws.Range("Y18").NumberFormat = "@"
$ws.Range("Y18").Value = "2016-06-02"

$ws.Range("X20").NumberFormat = "@"
$ws.Range("X20").Value = "2.2"
$ws.Range("Y20").NumberFormat = "@"
$ws.Range("Y20").Value = "2016-08-23"

$ws.Range("X21").NumberFormat = "@"
$ws.Range("X21").Value = "0.66"
$ws.Range("Y21").NumberFormat = "@"
$ws.Range("Y21").Value = "2015-09-21"

$ws.Range("X22").NumberFormat = "@"
$ws.Range("X22").Value = "3.23"
$ws.Range("Y22").NumberFormat = "@"
$ws.Range("Y22").Value = "2013-06-06"

$ws.Range("X23").NumberFormat = "@"
$ws.Range("X23").Value = "0.204"
$ws.Range("Y23").NumberFormat = "@"
$ws.Range("Y23").Value = "2017-11-01"

$ws.Range("X24").NumberFormat = "@"
$ws.Range("X24").Value = "2.6"
$ws.Range("Y24").NumberFormat = "@"
$ws.Range("Y24").Value = "2009-12-14"

$ws.Range("X25").NumberFormat = "@"
$ws.Range("X25").Value = "2.81"
$ws.Range("Y25").NumberFormat = "@"
$ws.Range("Y25").Value = "2014-04-09"

$ws.Range("X26").NumberFormat = "@"
$ws.Range("X26").Value = "0.843"
$ws.Range("Y26").NumberFormat = "@"
$ws.Range("Y26").Value = "2019-04-12"

$ws.Range("X27").NumberFormat = "@"
$ws.Range("X27").Value = "2.63"
$ws.Range("Y27").NumberFormat = "@"
$ws.Range("Y27").Value = "2013-12-11"
